# Straight from Karbon v1
# The sheet had a stale 45499-dated block (rows 2-10) and a correct
# 45474-dated block (rows 11-19). Replace the stale block's values with
# the correct ones (re-sorted by status: Clarification Required,
# Outreach Required, In Progress) and drop the now-redundant rows 11-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clarification Required
$ws.Range("A2").Value = "Clarification Required"
$ws.Range("B2").Value = "Low"
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 45474

$ws.Range("A3").Value = "Clarification Required"
$ws.Range("B3").Value = "Standard"
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 45474

$ws.Range("A4").Value = "Clarification Required"
$ws.Range("B4").Value = "Heightened"
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 45474

# Outreach Required
$ws.Range("A5").Value = "Outreach Required"
$ws.Range("B5").Value = "Low"
$ws.Range("C5").Value = 124
$ws.Range("D5").Value = 45474

$ws.Range("A6").Value = "Outreach Required"
$ws.Range("B6").Value = "Standard"
$ws.Range("C6").Value = 142
$ws.Range("D6").Value = 45474

$ws.Range("A7").Value = "Outreach Required"
$ws.Range("B7").Value = "Heightened"
$ws.Range("C7").Value = 91
$ws.Range("D7").Value = 45474

# In Progress
$ws.Range("A8").Value = "In Progress"
$ws.Range("B8").Value = "Low"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 45474

$ws.Range("A9").Value = "In Progress"
$ws.Range("B9").Value = "Standard"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 45474

$ws.Range("A10").Value = "In Progress"
$ws.Range("B10").Value = "Heightened"
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 45474

# The old rows 11-19 duplicated this data; remove them now that rows
# 2-10 hold the correct, re-sorted figures.
$ws.Range("A11:D19").EntireRow.Delete()
